$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 1")

$ws.Range("C9").Value = "Create User Stories"
$ws.Range("C13").Select()
